# Rearranged order of columns in test files.
# Swap column A and column D contents (columns B and C stay the same).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 3; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $dCell = $ws.Cells.Item($r, 4)

    $aVal = $aCell.Value()
    $dVal = $dCell.Value()

    $aCell.Value = $dVal
    $dCell.Value = $aVal
}

# Update the selection to reflect the new active range (A1:A3) with no
# specific active cell highlighted.
$ws.Range("A1:A3").Select()
